$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.656.41'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '2.374.67'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.677'
$ws.Range("E5").Value = '  +3.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.53'
$ws.Range("E6").Value = '  +3.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.49'
$ws.Range("E7").Value = '  +8.24%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  +19.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +9.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.41'
$ws.Range("E11").Value = '  +11.37%  '
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("D13").Value = '2.723.87'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.90'
$ws.Range("E14").Value = '  +9.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.72'
$ws.Range("E15").Value = '  +7.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.905'
$ws.Range("E16").Value = '  +7.56%  '
$ws.Range("D17").Value = '2.365.24'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '44.585.54'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("E19").Value = '  +6.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.03'
$ws.Range("E20").Value = '  +6.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.48'
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '255.90'
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  -4.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.52'
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.50'
$ws.Range("E26").Value = '  +5.14%  '
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.57'
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.01'
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("E30").Value = '  +5.81%  '
$ws.Range("E31").Value = '  +2.94%  '
$ws.Range("E32").Value = '  +5.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0741'
$ws.Range("E33").Value = '  +6.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.22'
$ws.Range("E34").Value = '  +4.68%  '
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.94'
$ws.Range("E36").Value = '  +8.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.45'
$ws.Range("E37").Value = '  -3.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.53'
$ws.Range("E38").Value = '  +0.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0273'
$ws.Range("E39").Value = '  +7.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.76'
$ws.Range("E40").Value = '  +10.52%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("E43").Value = '  +3.38%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.17'
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0982'
$ws.Range("E45").Value = '  +3.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.48'
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.78'
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.184'
$ws.Range("E48").Value = '  +12.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.36'
$ws.Range("E49").Value = '  +4.53%  '
$ws.Range("D50").Value = '1.443.04'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.78'
$ws.Range("E51").Value = '  +1.56%  '
